# Update the cryptos price list (column D = Price, column E = Volume(1h))
# with the latest scraped values.
#
# Note: a handful of Price values are numeric-looking text that ends in a
# significant trailing zero (e.g. "1.000", "155.90"). Assigning those bare
# to .Value would get auto-coerced to a number and silently drop the
# trailing zero, so they're written with a leading apostrophe, which is
# the standard Excel "treat as text" prefix and keeps the literal string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.928.02'
$ws.Cells.Item(2, 5).Value = '  +1.72%  '

$ws.Cells.Item(3, 4).Value = '1.890.88'
$ws.Cells.Item(3, 5).Value = '  +1.60%  '

$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  -0.49%  '

$ws.Cells.Item(5, 4).Value = '325.27'
$ws.Cells.Item(5, 5).Value = '  +0.04%  '

$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  -0.56%  '

$ws.Cells.Item(7, 4).Value = '0.4581'
$ws.Cells.Item(7, 5).Value = '  +0.68%  '

$ws.Cells.Item(8, 4).Value = '0.3902'
$ws.Cells.Item(8, 5).Value = '  +1.98%  '

$ws.Cells.Item(9, 4).Value = '0.07835'
$ws.Cells.Item(9, 5).Value = '  +0.33%  '

$ws.Cells.Item(10, 4).Value = '0.9883'
$ws.Cells.Item(10, 5).Value = '  +0.38%  '

$ws.Cells.Item(11, 4).Value = '21.86'
$ws.Cells.Item(11, 5).Value = '  +1.95%  '

$ws.Cells.Item(12, 4).Value = '1.955.38'
$ws.Cells.Item(12, 5).Value = '  +6.34%  '

$ws.Cells.Item(13, 4).Value = '7.026'
$ws.Cells.Item(13, 5).Value = '  +1.99%  '

$ws.Cells.Item(14, 4).Value = '5.682'
$ws.Cells.Item(14, 5).Value = '  +0.96%  '

$ws.Cells.Item(15, 4).Value = '0.06934'

$ws.Cells.Item(16, 4).Value = '87.99'
$ws.Cells.Item(16, 5).Value = '  +1.81%  '

$ws.Cells.Item(17, 5).Value = '  -0.63%  '

$ws.Cells.Item(18, 4).Value = '0.000009978'
$ws.Cells.Item(18, 5).Value = '  +0.60%  '

$ws.Cells.Item(19, 4).Value = '16.99'
$ws.Cells.Item(19, 5).Value = '  +2.05%  '

$ws.Cells.Item(20, 4).Value = '''1.000'
$ws.Cells.Item(20, 5).Value = '  -0.55%  '

$ws.Cells.Item(21, 4).Value = '28.919.84'
$ws.Cells.Item(21, 5).Value = '  +1.66%  '

$ws.Cells.Item(22, 4).Value = '5.295'
$ws.Cells.Item(22, 5).Value = '  +0.96%  '

$ws.Cells.Item(23, 5).Value = '  +1.09%  '

$ws.Cells.Item(24, 4).Value = '2.154.21'
$ws.Cells.Item(24, 5).Value = '  +4.39%  '

$ws.Cells.Item(25, 4).Value = '2.059'
$ws.Cells.Item(25, 5).Value = '  -1.40%  '

$ws.Cells.Item(26, 4).Value = '''155.90'
$ws.Cells.Item(26, 5).Value = '  +1.77%  '

$ws.Cells.Item(27, 4).Value = '19.26'
$ws.Cells.Item(27, 5).Value = '  +1.22%  '

$ws.Cells.Item(28, 4).Value = '5.892'
$ws.Cells.Item(28, 5).Value = '  +4.63%  '

$ws.Cells.Item(29, 4).Value = '1.925'
$ws.Cells.Item(29, 5).Value = '  +1.73%  '

$ws.Cells.Item(30, 4).Value = '117.45'
$ws.Cells.Item(30, 5).Value = '  +0.19%  '

$ws.Cells.Item(31, 5).Value = '  +0.83%  '

$ws.Cells.Item(32, 4).Value = '0.9049'
$ws.Cells.Item(32, 5).Value = '  +0.27%  '

$ws.Cells.Item(33, 5).Value = '  +0.51%  '

$ws.Cells.Item(34, 4).Value = '1.329'
$ws.Cells.Item(34, 5).Value = '  +1.18%  '

$ws.Cells.Item(35, 4).Value = '''3.260'
$ws.Cells.Item(35, 5).Value = '  -0.88%  '

$ws.Cells.Item(36, 4).Value = '1.188'
$ws.Cells.Item(36, 5).Value = '  +3.43%  '

$ws.Cells.Item(37, 4).Value = '''0.05770'
$ws.Cells.Item(37, 5).Value = '  +1.77%  '

$ws.Cells.Item(38, 4).Value = '''0.02070'
$ws.Cells.Item(38, 5).Value = '  +1.67%  '

$ws.Cells.Item(39, 4).Value = '''1.000'
$ws.Cells.Item(39, 5).Value = '  -0.62%  '

$ws.Cells.Item(40, 4).Value = '7.712'
$ws.Cells.Item(40, 5).Value = '  +1.19%  '

$ws.Cells.Item(41, 4).Value = '0.5679'
$ws.Cells.Item(41, 5).Value = '  +2.44%  '

$ws.Cells.Item(42, 4).Value = '''0.1770'
$ws.Cells.Item(42, 5).Value = '  +0.44%  '

$ws.Cells.Item(43, 4).Value = '9.727'

$ws.Cells.Item(44, 4).Value = '2.295'
$ws.Cells.Item(44, 5).Value = '  +9.00%  '

$ws.Cells.Item(45, 4).Value = '11.95'
$ws.Cells.Item(45, 5).Value = '  +4.17%  '

$ws.Cells.Item(46, 4).Value = '0.5344'
$ws.Cells.Item(46, 5).Value = '  +2.28%  '

$ws.Cells.Item(47, 4).Value = '0.07054'
$ws.Cells.Item(47, 5).Value = '  -0.94%  '

$ws.Cells.Item(48, 5).Value = '  +2.36%  '

$ws.Cells.Item(49, 4).Value = '112.84'
$ws.Cells.Item(49, 5).Value = '  +1.09%  '

$ws.Cells.Item(50, 4).Value = '2.526'
$ws.Cells.Item(50, 5).Value = '  +4.07%  '

$ws.Cells.Item(51, 4).Value = '''1.060'
$ws.Cells.Item(51, 5).Value = '  -5.67%  '
